$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 7.364689666666667
$ws.Range("H2").Value = 22.094069
$ws.Range("I2").Value = 0.5165094431700068
$ws.Range("J2").Value = 0.5165094431700068
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 44.69746666666666
$ws.Range("N2").Value = 134.0924
$ws.Range("O2").Value = 0.6823972194925493
$ws.Range("P2").Value = 0.6823972194925493
$ws.Range("Q2").Value = 329.1829708861778
$ws.Range("R2").Value = 2962.6467379756
$ws.Range("S2").Value = 0.3524646078608575
$ws.Range("T2").Value = 0.3524646078608575

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 7.364689666666667
$ws.Range("H3").Value = 22.094069
$ws.Range("I3").Value = 0.5165094431700068
$ws.Range("J3").Value = 0.5165094431700068
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 18.63243533333333
$ws.Range("N3").Value = 55.897306
$ws.Range("O3").Value = 0.2844618053784121
$ws.Range("P3").Value = 0.2844618053784121
$ws.Range("Q3").Value = 137.2221039642349
$ws.Range("R3").Value = 1234.998935678114
$ws.Range("S3").Value = 0.1469272086991385
$ws.Range("T3").Value = 0.1469272086991385

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 7.364689666666667
$ws.Range("H4").Value = 22.094069
$ws.Range("I4").Value = 0.5165094431700068
$ws.Range("J4").Value = 0.5165094431700068
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 2.170755666666667
$ws.Range("N4").Value = 6.512267
$ws.Range("O4").Value = 0.03314097512903853
$ws.Range("P4").Value = 0.03314097512903853
$ws.Range("Q4").Value = 15.98694182715811
$ws.Range("R4").Value = 143.882476444423
$ws.Range("S4").Value = 0.01711762661001074
$ws.Range("T4").Value = 0.01711762661001073

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 1.829962
$ws.Range("H5").Value = 5.489886
$ws.Range("I5").Value = 0.1283411381093639
$ws.Range("J5").Value = 0.1283411381093639
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 44.69746666666666
$ws.Range("N5").Value = 134.0924
$ws.Range("O5").Value = 0.6823972194925493
$ws.Range("P5").Value = 0.6823972194925493
$ws.Range("Q5").Value = 81.79466549626666
$ws.Range("R5").Value = 736.1519894664
$ws.Range("S5").Value = 0.0875796357923392
$ws.Range("T5").Value = 0.08757963579233918

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 1.829962
$ws.Range("H6").Value = 5.489886
$ws.Range("I6").Value = 0.1283411381093639
$ws.Range("J6").Value = 0.1283411381093639
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 18.63243533333333
$ws.Range("N6").Value = 55.897306
$ws.Range("O6").Value = 0.2844618053784121
$ws.Range("P6").Value = 0.2844618053784121
$ws.Range("Q6").Value = 34.09664862745733
$ws.Range("R6").Value = 306.869837647116
$ws.Range("S6").Value = 0.0365081518509098
$ws.Range("T6").Value = 0.03650815185090979

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 1.829962
$ws.Range("H7").Value = 5.489886
$ws.Range("I7").Value = 0.1283411381093639
$ws.Range("J7").Value = 0.1283411381093639
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 2.170755666666667
$ws.Range("N7").Value = 6.512267
$ws.Range("O7").Value = 0.03314097512903853
$ws.Range("P7").Value = 0.03314097512903853
$ws.Range("Q7").Value = 3.972400381284667
$ws.Range("R7").Value = 35.751603431562
$ws.Range("S7").Value = 0.004253350466114929
$ws.Range("T7").Value = 0.004253350466114927

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 5.063925333333334
$ws.Range("H8").Value = 15.191776
$ws.Range("I8").Value = 0.3551494187206292
$ws.Range("J8").Value = 0.3551494187206292
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 44.69746666666666
$ws.Range("N8").Value = 134.0924
$ws.Range("O8").Value = 0.6823972194925493
$ws.Range("P8").Value = 0.6823972194925493
$ws.Range("Q8").Value = 226.3446337891556
$ws.Range("R8").Value = 2037.1017041024
$ws.Range("S8").Value = 0.2423529758393525
$ws.Range("T8").Value = 0.2423529758393525

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 5.063925333333334
$ws.Range("H9").Value = 15.191776
$ws.Range("I9").Value = 0.3551494187206292
$ws.Range("J9").Value = 0.3551494187206292
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 18.63243533333333
$ws.Range("N9").Value = 55.897306
$ws.Range("O9").Value = 0.2844618053784121
$ws.Range("P9").Value = 0.2844618053784121
$ws.Range("Q9").Value = 94.35326130616178
$ws.Range("R9").Value = 849.1793517554561
$ws.Range("S9").Value = 0.1010264448283638
$ws.Range("T9").Value = 0.1010264448283638

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 5.063925333333334
$ws.Range("H10").Value = 15.191776
$ws.Range("I10").Value = 0.3551494187206292
$ws.Range("J10").Value = 0.3551494187206292
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 2.170755666666667
$ws.Range("N10").Value = 6.512267
$ws.Range("O10").Value = 0.03314097512903853
$ws.Range("P10").Value = 0.03314097512903853
$ws.Range("Q10").Value = 10.99254461291022
$ws.Range("R10").Value = 98.93290151619202
$ws.Range("S10").Value = 0.01176999805291287
$ws.Range("T10").Value = 0.01176999805291286
